$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "<__main__.Entity object at 0x000001E782C03090>"
$ws.Range("G2").Value = 4.558594024346748
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 7.549857829483546

# Row 3
$ws.Range("C3").Value = "<__main__.Entity object at 0x000001E7804AFF50>"
$ws.Range("G3").Value = 9.994117823433328
$ws.Range("I3").Value = 13.35831707259888

# Row 4
$ws.Range("C4").Value = "<__main__.Entity object at 0x000001E782DA9390>"
$ws.Range("G4").Value = 20.86020829870003
$ws.Range("I4").Value = 8.545549680531575

# Row 5
$ws.Range("C5").Value = "<__main__.Entity object at 0x000001E780309FD0>"
$ws.Range("G5").Value = 31.30957989156325
$ws.Range("I5").Value = 41.81996940997271

# Row 6
$ws.Range("C6").Value = "<__main__.Entity object at 0x000001E782CA2010>"
$ws.Range("G6").Value = 4.692896436942904
$ws.Range("I6").Value = 11.01765782958361

# Row 7
$ws.Range("C7").Value = "<__main__.Entity object at 0x000001E780221710>"
$ws.Range("G7").Value = 11.511948905071
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 9.092649555729189

# Row 8
$ws.Range("C8").Value = "<__main__.Entity object at 0x000001E780223550>"
$ws.Range("G8").Value = 4.449041314942946
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 4.152982905834425

# Row 9
$ws.Range("C9").Value = "<__main__.Entity object at 0x000001E780221750>"
$ws.Range("G9").Value = 14.58877561852719
$ws.Range("I9").Value = 23.69276483116812

# Row 10
$ws.Range("C10").Value = "<__main__.Entity object at 0x000001E78041D4D0>"
$ws.Range("G10").Value = 11.98991920880145
$ws.Range("I10").Value = 11.75015264334031

# Row 11
$ws.Range("C11").Value = "<__main__.Entity object at 0x000001E782A8FE50>"
$ws.Range("G11").Value = 6.647026366299166
$ws.Range("I11").Value = 4.131813628053326
